# Add a new "2022-Q1" worksheet (holdings detail) right before the "总计" (totals) sheet,
# and insert a corresponding summary row at the top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet positioned immediately before "总计"
# ---------------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheetBefore)
$q1.Name = "2022-Q1"

# NOTE: the "总计" sheet reference obtained before Worksheets.Add() becomes
# stale once the new sheet is inserted, so re-resolve it by name before any
# further use.
$totalSheet = $wb.Worksheets.Item("总计")

# Use an existing quarterly sheet as the source of header/row styling (bold
# header row + A-column index style) so the new sheet matches the look of
# its siblings.
$styleSource = $wb.Worksheets.Item("2021-Q4")

# Header row (row 1)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 2; $c -le 8; $c++) {
    $q1.Cells.Item(1, $c).Value = $headers[$c - 2]
}
$styleSource.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Holdings data rows (row 2 onward)
$rows = @(
    @("513060", "博时恒生医疗保健ETF（QDII）", "23.35", "98.89", "6.20", "1.4477", 4),
    @("011891", "易方达先锋成长混合A", "4.92", "94.38", "4.31", "0.2121", 10),
    @("011157", "弘毅远方港股通智选领航混合A", "3.65", "90.35", "5.19", "0.1894", 3),
    @("513700", "鹏华中证港股通医药卫生综合交易型开放式指数证券投资基金", "3.24", "93.11", "5.44", "0.1763", 5),
    @("517050", "华泰柏瑞中证沪港深互联网ETF", "6.18", "96.51", "2.83", "0.1749", 9),
    @("159792", "富国中证港股通互联网ETF", "2.76", "99.00", "5.75", "0.1587", 6),
    @("012379", "创金合信港股互联网3个月持有期混合型证券投资基金（QDII）A", "4.43", "73.95", "3.22", "0.1426", 10),
    @("159856", "工银瑞信中证沪港深互联网ETF", "4.44", "97.44", "3.06", "0.1359", 8),
    @("159892", "华夏恒生香港上市生物科技ETF（QDII）", "1.51", "99.03", "6.78", "0.1024", 4),
    @("517200", "嘉实中证沪港深互联网ETF", "1.71", "98.76", "3.12", "0.0534", 7),
    @("003993", "前海开源沪港深核心驱动灵活配置混合", "0.58", "82.10", "7.79", "0.0452", 7),
    @("159729", "汇添富中证沪港深互联网交易型开放式指数证券投资基金", "1.38", "97.96", "3.09", "0.0426", 8),
    @("011158", "弘毅远方港股通智选领航混合C", "0.78", "90.35", "5.19", "0.0405", 3),
    @("007151", "前海开源沪港深聚瑞混合", "0.60", "72.90", "6.53", "0.0392", 6),
    @("012380", "创金合信港股互联网3个月持有期混合型证券投资基金（QDII）C", "1.06", "73.95", "3.22", "0.0341", 10),
    @("012371", "西藏东财中证沪港深互联网指数型发起式证券投资基金A", "1.11", "95.04", "3.00", "0.0333", 7),
    @("012372", "西藏东财中证沪港深互联网指数型发起式证券投资基金C", "0.51", "95.04", "3.00", "0.0153", 7),
    @("011892", "易方达先锋成长混合C", "0.23", "94.38", "4.31", "0.0099", 10)
)

# Columns B-G must stay text (codes with leading zeros, and numbers whose
# trailing zeros matter, e.g. "6.20") so force text number-format first.
$lastRow = 1 + $rows.Count
$q1.Range("B2:G$lastRow").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $r - 2
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Apply the A-column index style (bold/centered, matching other sheets) to
# the newly written A2:A<lastRow> range.
$styleSource.Range("A2").Copy()
$q1.Range("A2:A$lastRow").PasteSpecial(-4122)

$q1.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of the "总计" sheet
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 18
$totalSheet.Range("D2").Value = 3.05

$totalSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3. Restore the originally active sheet/selection
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
